# Apply the changes described by the diff to the workbook.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "begroting"
$ws2 = $wb.Worksheets.Item(2)   # "Sheet1"

# --- Sheet "begroting" (sheet1) data changes -------------------------------

# Row 26: A26 becomes "week number" (a new shared string), and the existing
# "hours spent" label moves over to B26.
$ws1.Range("A26").Value = "week number"
$ws1.Range("B26").Value = "hours spent"

# B19 gets the built-in "Neutral" cell style applied (adds font/fill/style
# entries to styles.xml), text contents stay the same.
$ws1.Range("B19").Style = "Neutral"

# New rows with additional weekly hours data.
$ws1.Range("A38").Value = 15
$ws1.Range("B38").Value = 34
$ws1.Range("A39").Value = 16

# New "budget" label cell.
$ws1.Range("C41").Value = "budget"

# (B42 / C42 totals recalculate automatically via their existing formulas.)

# --- Sheet "Sheet1" (sheet2) data changes -----------------------------------

$ws2.Range("B6").Value = 5
# (B7 recalculates automatically via its existing formula B6/SQRT(9.81).)

# --- View / selection state --------------------------------------------------

# Before: "Sheet1" (sheet2) was the active/tab-selected sheet with selection
# N13, and "begroting" (sheet1) had selection B7:B21.
# After: "begroting" becomes the active/tab-selected sheet (so activeTab
# reverts to its default of 0) with selection B36, and "Sheet1" keeps a plain
# (non-active) view with selection B9.

$ws2.Range("B9").Select()

$ws1.Activate()
$ws1.Range("B36").Select()

Write-Host "edits applied"
